$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Save a reference plain (unstyled) style/number-format from a data cell to restore after forcing text format,
# so that writing number-like strings does not get auto-converted to numeric values nor leaves stray styling.
$plainStyle = $ws.Range("B2").Style

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "59.857.21"
$ws.Range("D2").Style = $plainStyle
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +1.64%  "
$ws.Range("E2").Style = $plainStyle

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.312.26"
$ws.Range("D3").Style = $plainStyle
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.09%  "
$ws.Range("E3").Style = $plainStyle

# Row 4
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("E4").Style = $plainStyle

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "540.52"
$ws.Range("D5").Style = $plainStyle
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.27%  "
$ws.Range("E5").Style = $plainStyle

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "129.84"
$ws.Range("D6").Style = $plainStyle
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -1.94%  "
$ws.Range("E6").Style = $plainStyle

# Row 7
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("E7").Style = $plainStyle

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.574"
$ws.Range("D8").Style = $plainStyle
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -2.13%  "
$ws.Range("E8").Style = $plainStyle

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.311.62"
$ws.Range("D9").Style = $plainStyle
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.17%  "
$ws.Range("E9").Style = $plainStyle

# Row 10
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -0.09%  "
$ws.Range("E10").Style = $plainStyle

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.48"
$ws.Range("D11").Style = $plainStyle
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +0.17%  "
$ws.Range("E11").Style = $plainStyle

# Row 12
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -0.13%  "
$ws.Range("E12").Style = $plainStyle

# Row 13
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -0.89%  "
$ws.Range("E13").Style = $plainStyle

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "23.33"
$ws.Range("D14").Style = $plainStyle
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -2.19%  "
$ws.Range("E14").Style = $plainStyle

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.725.50"
$ws.Range("D15").Style = $plainStyle
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +0.12%  "
$ws.Range("E15").Style = $plainStyle

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "59.869.65"
$ws.Range("D16").Style = $plainStyle
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +1.87%  "
$ws.Range("E16").Style = $plainStyle

# Row 17
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -1.23%  "
$ws.Range("E17").Style = $plainStyle

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.310.22"
$ws.Range("D18").Style = $plainStyle
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.41%  "
$ws.Range("E18").Style = $plainStyle

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.46"
$ws.Range("D19").Style = $plainStyle
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -1.90%  "
$ws.Range("E19").Style = $plainStyle

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.07"
$ws.Range("D20").Style = $plainStyle
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -2.25%  "
$ws.Range("E20").Style = $plainStyle

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "311.86"
$ws.Range("D21").Style = $plainStyle
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.26%  "
$ws.Range("E21").Style = $plainStyle

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.54"
$ws.Range("D22").Style = $plainStyle
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -0.99%  "
$ws.Range("E22").Style = $plainStyle

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.999"
$ws.Range("D23").Style = $plainStyle
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -0.26%  "
$ws.Range("E23").Style = $plainStyle

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.70"
$ws.Range("D24").Style = $plainStyle
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.19%  "
$ws.Range("E24").Style = $plainStyle

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "63.59"
$ws.Range("D25").Style = $plainStyle
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +2.02%  "
$ws.Range("E25").Style = $plainStyle

# Row 26
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -2.33%  "
$ws.Range("E26").Style = $plainStyle

# Row 27
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +0.01%  "
$ws.Range("E27").Style = $plainStyle

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.71"
$ws.Range("D28").Style = $plainStyle
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -3.00%  "
$ws.Range("E28").Style = $plainStyle

# Row 29
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +2.67%  "
$ws.Range("E29").Style = $plainStyle

# Row 30
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +0.14%  "
$ws.Range("E30").Style = $plainStyle

# Row 31
$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.70"
$ws.Range("D31").Style = $plainStyle
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -1.24%  "
$ws.Range("E31").Style = $plainStyle

# Row 32
$ws.Range("B32").Value = "SuiNetwork"
$ws.Range("C32").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.16"
$ws.Range("D32").Style = $plainStyle
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +0.43%  "
$ws.Range("E32").Style = $plainStyle

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0₃0722"
$ws.Range("D33").Style = $plainStyle
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -1.93%  "
$ws.Range("E33").Style = $plainStyle

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.80"
$ws.Range("D34").Style = $plainStyle
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -1.51%  "
$ws.Range("E34").Style = $plainStyle

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.35"
$ws.Range("D35").Style = $plainStyle
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +2.19%  "
$ws.Range("E35").Style = $plainStyle

# Row 36
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -2.16%  "
$ws.Range("E36").Style = $plainStyle

# Row 37
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +0.05%  "
$ws.Range("E37").Style = $plainStyle

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "17.66"
$ws.Range("D38").Style = $plainStyle
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -1.42%  "
$ws.Range("E38").Style = $plainStyle

# Row 39
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +0.20%  "
$ws.Range("E39").Style = $plainStyle

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.98"
$ws.Range("D40").Style = $plainStyle
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -2.71%  "
$ws.Range("E40").Style = $plainStyle

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "315.31"
$ws.Range("D41").Style = $plainStyle
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +4.38%  "
$ws.Range("E41").Style = $plainStyle

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "37.84"
$ws.Range("D42").Style = $plainStyle
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -1.49%  "
$ws.Range("E42").Style = $plainStyle

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.51"
$ws.Range("D43").Style = $plainStyle
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -0.40%  "
$ws.Range("E43").Style = $plainStyle

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "136.05"
$ws.Range("D44").Style = $plainStyle
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -4.01%  "
$ws.Range("E44").Style = $plainStyle

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.42"
$ws.Range("D45").Style = $plainStyle
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -1.01%  "
$ws.Range("E45").Style = $plainStyle

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0935"
$ws.Range("D46").Style = $plainStyle
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -2.54%  "
$ws.Range("E46").Style = $plainStyle

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.560"
$ws.Range("D47").Style = $plainStyle
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +0.72%  "
$ws.Range("E47").Style = $plainStyle

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "18.68"
$ws.Range("D48").Style = $plainStyle
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +1.93%  "
$ws.Range("E48").Style = $plainStyle

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0489"
$ws.Range("D49").Style = $plainStyle
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -1.46%  "
$ws.Range("E49").Style = $plainStyle

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0₆0220"
$ws.Range("D50").Style = $plainStyle
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +16.25%  "
$ws.Range("E50").Style = $plainStyle

# Row 51
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -0.35%  "
$ws.Range("E51").Style = $plainStyle
